# Streetbite deck - slide 7 "Streetbite Distance Algorithm"
# Reposition/resize the algorithm-listing placeholder and re-split its
# paragraphs into multiple runs (identifier tokens vs. plain text),
# matching the author's retype/spell-check pass. Overall text content
# is unchanged.

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(7)
$shp = $s.Shapes.Item(2)

# --- reposition / resize (values chosen so the EMU round-trip lands on
#     the exact target: 826618, 1166018, 8043062, 4579954) ---
$shp.Left   = 65.08803179606299
$shp.Top    = 91.81244284488189
$shp.Width  = 633.311981203937
$shp.Height = 360.62629921259844

$tr = $shp.TextFrame.TextRange

function Split-Paragraph {
    param([int]$Index, [string[]]$Segments)
    $para = $tr.Paragraphs($Index, 1)
    $pos = 1
    foreach ($seg in $Segments) {
        if ($seg.Length -gt 0) {
            $run = $para.Characters($pos, $seg.Length)
            $run.Text = $seg
        }
        $pos += $seg.Length
    }
}

Split-Paragraph 1 @("Algorithm ", "FindNearestVendor")
Split-Paragraph 2 @("Input: ", "userLat", ", ", "userLon")
Split-Paragraph 3 @("Output: Vendor closest to user")
Split-Paragraph 5 @("Step 1: Set ", "nearestVendor", " = null")
Split-Paragraph 6 @("Step 2: Set ", "smallestDistance", " = infinity")
Split-Paragraph 7 @("Step 3: For each vendor V:")
Split-Paragraph 8 @("           d = ", "ComputeDistance", "(user, V)")
Split-Paragraph 9 @("           If d < ", "smallestDistance", ":")
Split-Paragraph 10 @("               ", "smallestDistance", " = d")
Split-Paragraph 11 @("               ", "nearestVendor", " = V")
Split-Paragraph 12 @("Step 4: Return ", "nearestVendor")
Split-Paragraph 13 @("End Algorithm")
